$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "estoque_atualizado" (column G) values for the affected rows
$ws.Range("G2").Value = -154
$ws.Range("G4").Value = -22
$ws.Range("G5").Value = -148
$ws.Range("G11").Value = -26
